$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.123.89"
$ws.Range("E2").Value = "  -3.85%  "
$ws.Range("D3").Value = "2.588.78"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'566.44"
$ws.Range("E5").Value = "  -5.03%  "
$ws.Range("D6").Value = "'152.25"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("D9").Value = "2.588.51"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -8.87%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'0.372"
$ws.Range("E13").Value = "  -6.32%  "
$ws.Range("D14").Value = "'27.65"
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("D15").Value = "3.052.90"
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("E16").Value = "  -9.34%  "
$ws.Range("D17").Value = "63.002.84"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "2.594.36"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "'11.79"
$ws.Range("E19").Value = "  -5.43%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  -7.23%  "
$ws.Range("D22").Value = "'336.13"
$ws.Range("E22").Value = "  -4.73%  "
$ws.Range("D24").Value = "'66.64"
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -7.46%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'573.46"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("D29").Value = "'1.52"
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("E32").Value = "  -5.49%  "
$ws.Range("D33").Value = "'2.01"
$ws.Range("E33").Value = "  -5.38%  "
$ws.Range("D34").Value = "'1.67"
$ws.Range("E34").Value = "  -7.02%  "
$ws.Range("D35").Value = "'6.40"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("D39").Value = "'19.41"
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("D40").Value = "'154.37"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -6.66%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'41.25"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'156.17"
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("D46").Value = "'22.85"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'3.79"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("E48").Value = "  -6.54%  "
$ws.Range("D49").Value = "'0.622"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'0.0984"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("E51").Value = "  -5.84%  "
